# Auto-update draw results: append the new Pick 4 draw row (2025-10-12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26
$newRange = $ws.Range("A26:E26")

# The sheet stores every value (dates, numeric-looking phase/result codes,
# timestamps) as literal text. Force text interpretation first so Excel's
# auto-conversion doesn't turn "2025-10-12" into a date serial or
# "251012" into a number, then restore the default "Normal" style so the
# new row doesn't pick up a stray explicit style compared to its neighbours.
$newRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-10-12"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "251012"
$ws.Cells.Item($row, 4).Value = "4-9-2-1"
$ws.Cells.Item($row, 5).Value = "2025-10-12T21:35:22.246+04:00"

$newRange.Style = "Normal"
